$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value, and whether the text
# "looks like a number" to Excel's auto-detection (so it must be forced
# to remain Text, matching the original inlineStr cell type).
$updates = @(
    @{ Cell = 'D2'; Value = '25.861.16'; Numeric = $false }
    @{ Cell = 'E2'; Value = '  -0.15%  '; Numeric = $false }
    @{ Cell = 'D3'; Value = '1.631.12'; Numeric = $false }
    @{ Cell = 'E3'; Value = '  -0.58%  '; Numeric = $false }
    @{ Cell = 'D4'; Value = '1.003'; Numeric = $true }
    @{ Cell = 'E4'; Value = '  +0.23%  '; Numeric = $false }
    @{ Cell = 'D5'; Value = '213.97'; Numeric = $true }
    @{ Cell = 'E5'; Value = '  -0.48%  '; Numeric = $false }
    @{ Cell = 'D6'; Value = '0.5094'; Numeric = $true }
    @{ Cell = 'E6'; Value = '  +1.02%  '; Numeric = $false }
    @{ Cell = 'D7'; Value = '1.002'; Numeric = $true }
    @{ Cell = 'E7'; Value = '  -0.25%  '; Numeric = $false }
    @{ Cell = 'D8'; Value = '0.2546'; Numeric = $true }
    @{ Cell = 'E8'; Value = '  -1.08%  '; Numeric = $false }
    @{ Cell = 'D9'; Value = '0.06324'; Numeric = $true }
    @{ Cell = 'E9'; Value = '  -1.07%  '; Numeric = $false }
    @{ Cell = 'E10'; Value = '  -0.65%  '; Numeric = $false }
    @{ Cell = 'D11'; Value = '0.07734'; Numeric = $true }
    @{ Cell = 'E11'; Value = '  -0.46%  '; Numeric = $false }
    @{ Cell = 'D12'; Value = '4.267'; Numeric = $true }
    @{ Cell = 'E12'; Value = '  -0.02%  '; Numeric = $false }
    @{ Cell = 'D13'; Value = '1.632.16'; Numeric = $false }
    @{ Cell = 'E13'; Value = '  -0.66%  '; Numeric = $false }
    @{ Cell = 'D14'; Value = '0.5403'; Numeric = $true }
    @{ Cell = 'E14'; Value = '  -0.56%  '; Numeric = $false }
    @{ Cell = 'D15'; Value = '0.0₅7687'; Numeric = $false }
    @{ Cell = 'E15'; Value = '  -2.88%  '; Numeric = $false }
    @{ Cell = 'D16'; Value = '63.95'; Numeric = $true }
    @{ Cell = 'E16'; Value = '  -0.93%  '; Numeric = $false }
    @{ Cell = 'D17'; Value = '25.871.91'; Numeric = $false }
    @{ Cell = 'E18'; Value = '  -0.43%  '; Numeric = $false }
    @{ Cell = 'D19'; Value = '194.46'; Numeric = $true }
    @{ Cell = 'E19'; Value = '  -1.49%  '; Numeric = $false }
    @{ Cell = 'D20'; Value = '4.408'; Numeric = $true }
    @{ Cell = 'E20'; Value = '  +0.82%  '; Numeric = $false }
    @{ Cell = 'D21'; Value = '9.868'; Numeric = $true }
    @{ Cell = 'E21'; Value = '  -0.35%  '; Numeric = $false }
    @{ Cell = 'D22'; Value = '5.997'; Numeric = $true }
    @{ Cell = 'E22'; Value = '  +0.58%  '; Numeric = $false }
    @{ Cell = 'E23'; Value = '  -0.33%  '; Numeric = $false }
    @{ Cell = 'D24'; Value = '1.860'; Numeric = $true }
    @{ Cell = 'E24'; Value = '  -0.69%  '; Numeric = $false }
    @{ Cell = 'D25'; Value = '140.58'; Numeric = $true }
    @{ Cell = 'E25'; Value = '  -0.54%  '; Numeric = $false }
    @{ Cell = 'E26'; Value = '  +4.15%  '; Numeric = $false }
    @{ Cell = 'E27'; Value = '  -0.57%  '; Numeric = $false }
    @{ Cell = 'D28'; Value = '15.49'; Numeric = $true }
    @{ Cell = 'E28'; Value = '  -1.13%  '; Numeric = $false }
    @{ Cell = 'E29'; Value = '  -0.64%  '; Numeric = $false }
    @{ Cell = 'D30'; Value = '0.04890'; Numeric = $true }
    @{ Cell = 'E30'; Value = '  -1.02%  '; Numeric = $false }
    @{ Cell = 'D31'; Value = '3.232'; Numeric = $true }
    @{ Cell = 'E31'; Value = '  -0.95%  '; Numeric = $false }
    @{ Cell = 'D32'; Value = '3.153'; Numeric = $true }
    @{ Cell = 'E32'; Value = '  -1.37%  '; Numeric = $false }
    @{ Cell = 'D33'; Value = '1.519'; Numeric = $true }
    @{ Cell = 'E33'; Value = '  -1.30%  '; Numeric = $false }
    @{ Cell = 'D34'; Value = '2.363'; Numeric = $true }
    @{ Cell = 'E34'; Value = '  -0.58%  '; Numeric = $false }
    @{ Cell = 'D35'; Value = '0.8855'; Numeric = $true }
    @{ Cell = 'D36'; Value = '2.572'; Numeric = $true }
    @{ Cell = 'E36'; Value = '  -1.41%  '; Numeric = $false }
    @{ Cell = 'D37'; Value = '1.136.05'; Numeric = $false }
    @{ Cell = 'E37'; Value = '  -0.82%  '; Numeric = $false }
    @{ Cell = 'D38'; Value = '0.5380'; Numeric = $true }
    @{ Cell = 'E38'; Value = '  -2.97%  '; Numeric = $false }
    @{ Cell = 'D39'; Value = '0.01543'; Numeric = $true }
    @{ Cell = 'E39'; Value = '  -1.74%  '; Numeric = $false }
    @{ Cell = 'E40'; Value = '  -0.58%  '; Numeric = $false }
    @{ Cell = 'D41'; Value = '2.531'; Numeric = $true }
    @{ Cell = 'E41'; Value = '  -1.08%  '; Numeric = $false }
    @{ Cell = 'E42'; Value = '  +4.82%  '; Numeric = $false }
    @{ Cell = 'D43'; Value = '0.8105'; Numeric = $true }
    @{ Cell = 'E43'; Value = '  -0.23%  '; Numeric = $false }
    @{ Cell = 'D44'; Value = '98.47'; Numeric = $true }
    @{ Cell = 'E44'; Value = '  -1.36%  '; Numeric = $false }
    @{ Cell = 'D45'; Value = '5.421'; Numeric = $true }
    @{ Cell = 'E45'; Value = '  -4.97%  '; Numeric = $false }
    @{ Cell = 'D46'; Value = '1.767.42'; Numeric = $false }
    @{ Cell = 'E46'; Value = '  -0.61%  '; Numeric = $false }
    @{ Cell = 'D47'; Value = '0.4530'; Numeric = $true }
    @{ Cell = 'E47'; Value = '  +0.46%  '; Numeric = $false }
    @{ Cell = 'E48'; Value = '  +0.10%  '; Numeric = $false }
    @{ Cell = 'D49'; Value = '54.49'; Numeric = $true }
    @{ Cell = 'E49'; Value = '  -0.51%  '; Numeric = $false }
    @{ Cell = 'D50'; Value = '0.05050'; Numeric = $true }
    @{ Cell = 'E51'; Value = '  -0.41%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # Force text storage so e.g. "1.003" / "0.5094" are not
        # reinterpreted as numbers, then restore the default "Normal"
        # style so no stray number-format style gets attached to the cell.
        $cell.NumberFormat = '@'
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
